$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2..16 (Sending cluster x Target cluster cross join,
# now including the new "M2" target cluster alongside ECs, FAPs, M1, sCs).
$data = @(
    @("ECs", "Sema6a", "Plxna4", "ECs", 3, 1, 51.12915033333334, 153.387451, 0.8013178159252168, 0.8013178159252169, 3, 1, 5.553617333333334, 16.660852, 0.6738463085322981, 0.6738463085322983, 283.9517355298058, 2555.565619768252, 0.5399650522223709, 0.5399650522223711),
    @("ECs", "Sema6a", "Plxna4", "FAPs", 3, 1, 51.12915033333334, 153.387451, 0.8013178159252168, 0.8013178159252169, 3, 1, 2.371913666666666, 7.115741, 0.2877953543625453, 0.2877953543625454, 121.2739304406879, 1091.465373966191, 0.2306155447912186, 0.2306155447912188),
    @("ECs", "Sema6a", "Plxna4", "M1", 3, 1, 51.12915033333334, 153.387451, 0.8013178159252168, 0.8013178159252169, 3, 1, 0.075423, 0.226269, 0.009151424572122395, 0.009151424572122398, 3.856313905591001, 34.706825150319, 0.007333199550737479, 0.007333199550737483),
    @("ECs", "Sema6a", "Plxna4", "M2", 3, 1, 51.12915033333334, 153.387451, 0.8013178159252168, 0.8013178159252169, 1, 0.3333333333333333, 0.11184, 0.33552, 0.01357006913204418, 0.01357006913204419, 5.71828417328, 51.46455755952, 0.01087393815884385, 0.01087393815884385),
    @("ECs", "Sema6a", "Plxna4", "sCs", 3, 1, 51.12915033333334, 153.387451, 0.8013178159252168, 0.8013178159252169, 3, 1, 0.1288736666666667, 0.386621, 0.01563684340098967, 0.01563684340098968, 6.589201077007889, 59.302809693071, 0.01253008120204568, 0.01253008120204569),
    @("FAPs", "Sema6a", "Plxna4", "ECs", 3, 1, 2.875170333333333, 8.625510999999999, 0.04506089377389114, 0.04506089377389114, 3, 1, 5.553617333333334, 16.660852, 0.6738463085322981, 0.6738463085322983, 15.96759579948578, 143.708362195372, 0.03036411692870256, 0.03036411692870257),
    @("FAPs", "Sema6a", "Plxna4", "FAPs", 3, 1, 2.875170333333333, 8.625510999999999, 0.04506089377389114, 0.04506089377389114, 3, 1, 2.371913666666666, 7.115741, 0.2877953543625453, 0.2877953543625454, 6.819655807627888, 61.37690226865099, 0.01296831589155001, 0.01296831589155002),
    @("FAPs", "Sema6a", "Plxna4", "M1", 3, 1, 2.875170333333333, 8.625510999999999, 0.04506089377389114, 0.04506089377389114, 3, 1, 0.075423, 0.226269, 0.009151424572122395, 0.009151424572122398, 0.216853972051, 1.951685748459, 0.0004123713705241844, 0.0004123713705241846),
    @("FAPs", "Sema6a", "Plxna4", "M2", 3, 1, 2.875170333333333, 8.625510999999999, 0.04506089377389114, 0.04506089377389114, 1, 0.3333333333333333, 0.11184, 0.33552, 0.01357006913204418, 0.01357006913204419, 0.32155905008, 2.89403145072, 0.0006114794436634021, 0.0006114794436634024),
    @("FAPs", "Sema6a", "Plxna4", "sCs", 3, 1, 2.875170333333333, 8.625510999999999, 0.04506089377389114, 0.04506089377389114, 3, 1, 0.1288736666666667, 0.386621, 0.01563684340098967, 0.01563684340098968, 0.3705337431478888, 3.334803688331, 0.0007046101394509662, 0.0007046101394509666),
    @("sCs", "Sema6a", "Plxna4", "ECs", 3, 1, 9.802011, 29.406033, 0.153621290300892, 0.153621290300892, 3, 1, 5.553617333333334, 16.660852, 0.6738463085322981, 0.6738463085322983, 54.43661819112401, 489.9295637201161, 0.1035171393812246, 0.1035171393812247),
    @("sCs", "Sema6a", "Plxna4", "FAPs", 3, 1, 9.802011, 29.406033, 0.153621290300892, 0.153621290300892, 3, 1, 2.371913666666666, 7.115741, 0.2877953543625453, 0.2877953543625454, 23.249523851717, 209.245714665453, 0.04421149367977667, 0.04421149367977668),
    @("sCs", "Sema6a", "Plxna4", "M1", 3, 1, 9.802011, 29.406033, 0.153621290300892, 0.153621290300892, 3, 1, 0.075423, 0.226269, 0.009151424572122395, 0.009151424572122398, 0.7392970756530001, 6.653673680877001, 0.001405853650860731, 0.001405853650860732),
    @("sCs", "Sema6a", "Plxna4", "M2", 3, 1, 9.802011, 29.406033, 0.153621290300892, 0.153621290300892, 1, 0.3333333333333333, 0.11184, 0.33552, 0.01357006913204418, 0.01357006913204419, 1.09625691024, 9.866312192160001, 0.002084651529536933, 0.002084651529536934),
    @("sCs", "Sema6a", "Plxna4", "sCs", 3, 1, 9.802011, 29.406033, 0.153621290300892, 0.153621290300892, 3, 1, 0.1288736666666667, 0.386621, 0.01563684340098967, 0.01563684340098968, 1.263221098277, 11.368989884493, 0.002402152059493022, 0.002402152059493023)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
